$wb = $excel.ActiveWorkbook

# --- ALC (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5396.7144
$ws.Range("J17").Value = 5875.316
$ws.Range("L17").Value = 17625.948
$ws.Range("N17").Value = -17961.948

# --- ALC (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4997
$ws.Range("J64").Value = 4997
$ws.Range("L64").Value = 4997
$ws.Range("N64").Value = -5493

# --- ALC (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4997
$ws.Range("J67").Value = 4997
$ws.Range("L67").Value = 4997
$ws.Range("N67").Value = -6713

# --- ALC (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4555.6875
$ws.Range("I86").Value = 3421.5454
$ws.Range("K86").Value = 3421.5454
$ws.Range("M86").Value = -2298.5454

# --- ALC (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4555.6875
$ws.Range("I89").Value = 3421.5454
$ws.Range("K89").Value = 17107.727
$ws.Range("M89").Value = -11491.727

# --- ALC (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1400.5714
$ws.Range("I98").Value = 1400.5714
$ws.Range("K98").Value = 1400.5714
$ws.Range("M98").Value = 97.42859999999996

# --- ALC (hunk 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8568.261
$ws.Range("I106").Value = 3378.8333
$ws.Range("J106").Value = 14229.454
$ws.Range("K106").Value = 3378.8333
$ws.Range("L106").Value = 14229.454
$ws.Range("M106").Value = -2747.8333
$ws.Range("N106").Value = -15491.454

# --- ALC (hunk 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1688.8334
$ws.Range("J112").Value = 2754.7144
$ws.Range("L112").Value = 8264.143199999999
$ws.Range("N112").Value = -10480.1432

# --- ALC (hunk 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1400.5714
$ws.Range("I122").Value = 1400.5714
$ws.Range("K122").Value = 4201.7142
$ws.Range("M122").Value = -1751.7142

# --- ARM (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 66670836
$ws.Range("J45").Value = 6876.5
$ws.Range("L45").Value = 6876.5
$ws.Range("N45").Value = -7630.5

# --- ARM (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6062.6
$ws.Range("I132").Value = 2164.6667
$ws.Range("K132").Value = 6494.000100000001
$ws.Range("M132").Value = -3964.000100000001

# --- BSM (hunk 11) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1875.3636
$ws.Range("J64").Value = 2358.3333
$ws.Range("L64").Value = 2358.3333
$ws.Range("N64").Value = -2808.3333

# --- BSM (hunk 12) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1875.3636
$ws.Range("J67").Value = 2358.3333
$ws.Range("L67").Value = 2358.3333
$ws.Range("N67").Value = -3918.3333

# --- BSM (hunk 13) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1100.1818
$ws.Range("I94").Value = 1130.2
$ws.Range("K94").Value = 1130.2
$ws.Range("M94").Value = -679.2

# --- BSM (hunk 14) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 64500
$ws.Range("J133").Value = 64500
$ws.Range("L133").Value = 64500
$ws.Range("N133").Value = -74620

# --- CRP (hunk 15) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258
$ws.Range("I7").Value = 209.55556
$ws.Range("J7").Value = 330.66666
$ws.Range("K7").Value = 209.55556
$ws.Range("L7").Value = 330.66666
$ws.Range("M7").Value = -96.55556000000001
$ws.Range("N7").Value = -556.66666

# --- CRP (hunk 16) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33459.473
$ws.Range("I31").Value = 3182.1365
$ws.Range("K31").Value = 3182.1365
$ws.Range("M31").Value = -2887.1365

# --- CRP (hunk 17) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 33459.473
$ws.Range("I34").Value = 3182.1365
$ws.Range("K34").Value = 3182.1365
$ws.Range("M34").Value = -2980.1365

# --- CRP (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4455.457
$ws.Range("I58").Value = 3510.889
$ws.Range("J58").Value = 7643.375
$ws.Range("K58").Value = 3510.889
$ws.Range("L58").Value = 7643.375
$ws.Range("M58").Value = -3307.889
$ws.Range("N58").Value = -8049.375

# --- CRP (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3900.5435
$ws.Range("I132").Value = 3675.2559
$ws.Range("K132").Value = 11025.7677
$ws.Range("M132").Value = -8495.7677

# --- CRP (hunk 20) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 61662.332
$ws.Range("J133").Value = 62495
$ws.Range("L133").Value = 62495
$ws.Range("N133").Value = -67555

# --- CRP (hunk 21) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4455.457
$ws.Range("I136").Value = 3510.889
$ws.Range("J136").Value = 7643.375
$ws.Range("K136").Value = 10532.667
$ws.Range("L136").Value = 22930.125
$ws.Range("M136").Value = -7982.667000000001
$ws.Range("N136").Value = -28030.125

# --- CUL (hunk 22) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 50182.2
$ws.Range("J7").Value = 83532
$ws.Range("L7").Value = 250596
$ws.Range("N7").Value = -250820

# --- CUL (hunk 23) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2181.9092
$ws.Range("I34").Value = 460.4
$ws.Range("J34").Value = 3616.5
$ws.Range("K34").Value = 1381.2
$ws.Range("L34").Value = 10849.5
$ws.Range("M34").Value = -1297.2
$ws.Range("N34").Value = -11017.5

# --- CUL (hunk 24) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3608.818
$ws.Range("I39").Value = 3499.6667
$ws.Range("J39").Value = 3649.75
$ws.Range("K39").Value = 10499.0001
$ws.Range("L39").Value = 10949.25
$ws.Range("M39").Value = -10205.0001
$ws.Range("N39").Value = -11537.25

# --- CUL (hunk 25) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2984.6
$ws.Range("I55").Value = 1108
$ws.Range("J55").Value = 5799.5
$ws.Range("K55").Value = 3324
$ws.Range("L55").Value = 17398.5
$ws.Range("M55").Value = -3147
$ws.Range("N55").Value = -17752.5

# --- GSM (hunk 26) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3920
$ws.Range("I102").Value = 2145.3333
$ws.Range("K102").Value = 2145.3333
$ws.Range("M102").Value = -523.3332999999998

# --- GSM (hunk 27) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6974.2856
$ws.Range("I122").Value = 3602.3333
$ws.Range("K122").Value = 10806.9999
$ws.Range("M122").Value = -8356.999899999999

# --- LTW (hunk 28) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3161.1875
$ws.Range("I16").Value = 2506.2307
$ws.Range("J16").Value = 5999.3335
$ws.Range("K16").Value = 2506.2307
$ws.Range("L16").Value = 5999.3335
$ws.Range("M16").Value = -2336.2307
$ws.Range("N16").Value = -6339.3335

# --- LTW (hunk 29) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7408.4
$ws.Range("I22").Value = 1980.3334
$ws.Range("K22").Value = 1980.3334
$ws.Range("M22").Value = -1685.3334

# --- LTW (hunk 30) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7408.4
$ws.Range("I27").Value = 1980.3334
$ws.Range("K27").Value = 1980.3334
$ws.Range("M27").Value = -1873.3334

# --- LTW (hunk 31) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3494.16
$ws.Range("I82").Value = 1041.7646
$ws.Range("J82").Value = 8705.5
$ws.Range("K82").Value = 1041.7646
$ws.Range("L82").Value = 8705.5
$ws.Range("M82").Value = -680.7646
$ws.Range("N82").Value = -9427.5

# --- LTW (hunk 32) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3494.16
$ws.Range("I85").Value = 1041.7646
$ws.Range("J85").Value = 8705.5
$ws.Range("K85").Value = 1041.7646
$ws.Range("L85").Value = 8705.5
$ws.Range("M85").Value = 206.2354
$ws.Range("N85").Value = -11201.5

# --- LTW (hunk 33) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4756.278
$ws.Range("I132").Value = 3064.3635
$ws.Range("K132").Value = 9193.0905
$ws.Range("M132").Value = -6663.0905

# --- WVR (hunk 34) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6000.3335
$ws.Range("I62").Value = 6000.3335
$ws.Range("K62").Value = 6000.3335
$ws.Range("M62").Value = -5376.3335

# --- WVR (hunk 35) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6000.3335
$ws.Range("I65").Value = 6000.3335
$ws.Range("K65").Value = 30001.6675
$ws.Range("M65").Value = -26881.6675

# --- WVR (hunk 36) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 49999
$ws.Range("J75").Value = 49999
$ws.Range("L75").Value = 49999
$ws.Range("N75").Value = -51871

# --- WVR (hunk 37) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 49999
$ws.Range("J78").Value = 49999
$ws.Range("L78").Value = 149997
$ws.Range("N78").Value = -159357

# --- WVR (hunk 38) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1409.6666
$ws.Range("I100").Value = 596.8333
$ws.Range("K100").Value = 1193.6666
$ws.Range("M100").Value = -652.6666

# --- WVR (hunk 39) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 748
$ws.Range("I107").Value = 999.6667
$ws.Range("K107").Value = 2999.0001
$ws.Range("M107").Value = -1079.0001

# --- WVR (hunk 40) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
